# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.17 = 12174.33 pesos`n✅ 12174.33 pesos = 3.16 = 958.63 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update numeric rate values on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("O10").Value = 3841
$wsTasas.Range("N12").Value = 3848
$wsTasas.Range("O12").Value = 303
